$wb = $excel.ActiveWorkbook

# --- Sheet ALC: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1436.5454
$ws.Range("I98").Value = 518.7368
$ws.Range("J98").Value = 7249.3335
$ws.Range("K98").Value = 518.7368
$ws.Range("L98").Value = 7249.3335
$ws.Range("M98").Value = 979.2632
$ws.Range("N98").Value = -10245.3335
$ws.Range("H112").Value = 1436.4706
$ws.Range("J112").Value = 1608.9286
$ws.Range("L112").Value = 4826.7858
$ws.Range("N112").Value = -7042.7858
$ws.Range("H122").Value = 1436.5454
$ws.Range("I122").Value = 518.7368
$ws.Range("J122").Value = 7249.3335
$ws.Range("K122").Value = 1556.2104
$ws.Range("L122").Value = 21748.0005
$ws.Range("M122").Value = 893.7896000000001
$ws.Range("N122").Value = -26648.0005
$ws.Range("H132").Value = 1996.3334
$ws.Range("I132").Value = 1321.62
$ws.Range("K132").Value = 3964.86
$ws.Range("M132").Value = -1434.86

# --- Sheet ARM: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2777.7778
$ws.Range("I102").Value = 1714.2858
$ws.Range("J102").Value = 6500
$ws.Range("K102").Value = 1714.2858
$ws.Range("L102").Value = 6500
$ws.Range("M102").Value = -92.28580000000011
$ws.Range("N102").Value = -9744

# --- Sheet BSM: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4341.905
$ws.Range("I86").Value = 9733.333000000001
$ws.Range("J86").Value = 2185.3333
$ws.Range("K86").Value = 9733.333000000001
$ws.Range("L86").Value = 2185.3333
$ws.Range("M86").Value = -8610.333000000001
$ws.Range("N86").Value = -4431.3333
$ws.Range("H89").Value = 4341.905
$ws.Range("I89").Value = 9733.333000000001
$ws.Range("J89").Value = 2185.3333
$ws.Range("K89").Value = 48666.665
$ws.Range("L89").Value = 10926.6665
$ws.Range("M89").Value = -43050.665
$ws.Range("N89").Value = -22158.6665

# --- Sheet CRP: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1484.5172
$ws.Range("I58").Value = 861.1667
$ws.Range("J58").Value = 2504.5454
$ws.Range("K58").Value = 861.1667
$ws.Range("L58").Value = 2504.5454
$ws.Range("M58").Value = -658.1667
$ws.Range("N58").Value = -2910.5454
$ws.Range("H136").Value = 1484.5172
$ws.Range("I136").Value = 861.1667
$ws.Range("J136").Value = 2504.5454
$ws.Range("K136").Value = 2583.5001
$ws.Range("L136").Value = 7513.6362
$ws.Range("M136").Value = -33.5001000000002
$ws.Range("N136").Value = -12613.6362

# --- Sheet CUL: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4333.4287
$ws.Range("I3").Value = 3555.6667
$ws.Range("K3").Value = 10667.0001
$ws.Range("M3").Value = -10555.0001
$ws.Range("H18").Value = 3161.9092
$ws.Range("I18").Value = 3683.111
$ws.Range("K18").Value = 11049.333
$ws.Range("M18").Value = -10880.333
$ws.Range("H62").Value = 2974.5
$ws.Range("I62").Value = 2900
$ws.Range("J62").Value = 2999.3333
$ws.Range("K62").Value = 8700
$ws.Range("L62").Value = 8997.999899999999
$ws.Range("M62").Value = -8014
$ws.Range("N62").Value = -10369.9999
$ws.Range("H65").Value = 2974.5
$ws.Range("I65").Value = 2900
$ws.Range("J65").Value = 2999.3333
$ws.Range("K65").Value = 26100
$ws.Range("L65").Value = 26993.9997
$ws.Range("M65").Value = -22668
$ws.Range("N65").Value = -33857.9997
$ws.Range("H81").Value = 2779
$ws.Range("I81").Value = 1202.6
$ws.Range("J81").Value = 4749.5
$ws.Range("K81").Value = 3607.8
$ws.Range("L81").Value = 14248.5
$ws.Range("M81").Value = -2484.8
$ws.Range("N81").Value = -16494.5
$ws.Range("H84").Value = 2779
$ws.Range("I84").Value = 1202.6
$ws.Range("J84").Value = 4749.5
$ws.Range("K84").Value = 10823.4
$ws.Range("L84").Value = 42745.5
$ws.Range("M84").Value = -5207.4
$ws.Range("N84").Value = -53977.5
$ws.Range("H113").Value = 13514054
$ws.Range("I113").Value = 19231322
$ws.Range("J113").Value = 510
$ws.Range("K113").Value = 57693966
$ws.Range("L113").Value = 1530
$ws.Range("M113").Value = -57691796
$ws.Range("N113").Value = -5870
$ws.Range("H122").Value = 25000598
$ws.Range("I122").Value = 38461916
$ws.Range("J122").Value = 1007.1429
$ws.Range("K122").Value = 346157244
$ws.Range("L122").Value = 9064.286100000001
$ws.Range("M122").Value = -346154794
$ws.Range("N122").Value = -13964.2861
$ws.Range("H131").Value = 1850.8823
$ws.Range("I131").Value = 2864.1667
$ws.Range("J131").Value = 1539.1025
$ws.Range("K131").Value = 8592.500100000001
$ws.Range("L131").Value = 4617.3075
$ws.Range("M131").Value = -3552.500100000001
$ws.Range("N131").Value = -14697.3075
$ws.Range("H133").Value = 3623.75
$ws.Range("I133").Value = 1548.3334
$ws.Range("J133").Value = 9850
$ws.Range("K133").Value = 4645.0002
$ws.Range("L133").Value = 29550
$ws.Range("M133").Value = 414.9997999999996
$ws.Range("N133").Value = -39670

# --- Sheet GSM: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3249.2969
$ws.Range("I80").Value = 3364.0566
$ws.Range("J80").Value = 2696.3635
$ws.Range("K80").Value = 3364.0566
$ws.Range("L80").Value = 2696.3635
$ws.Range("M80").Value = -2366.0566
$ws.Range("N80").Value = -4692.363499999999
$ws.Range("H83").Value = 3249.2969
$ws.Range("I83").Value = 3364.0566
$ws.Range("J83").Value = 2696.3635
$ws.Range("K83").Value = 16820.283
$ws.Range("L83").Value = 13481.8175
$ws.Range("M83").Value = -11828.283
$ws.Range("N83").Value = -23465.8175

# --- Sheet LTW: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 32749.834
$ws.Range("I88").Value = 9000
$ws.Range("J88").Value = 37499.8
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 37499.8
$ws.Range("M88").Value = -8572
$ws.Range("N88").Value = -38355.8
$ws.Range("H91").Value = 32749.834
$ws.Range("I91").Value = 9000
$ws.Range("J91").Value = 37499.8
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 37499.8
$ws.Range("M91").Value = -7518
$ws.Range("N91").Value = -40463.8

# --- Sheet WVR: update LeveProfit/price calculations (scheduled runner refresh) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29913.334
$ws.Range("J63").Value = 29913.334
$ws.Range("L63").Value = 29913.334
$ws.Range("N63").Value = -31161.334
$ws.Range("H66").Value = 29913.334
$ws.Range("J66").Value = 29913.334
$ws.Range("L66").Value = 89740.00199999999
$ws.Range("N66").Value = -95980.00199999999
$ws.Range("H69").Value = 32710
$ws.Range("J69").Value = 32710
$ws.Range("L69").Value = 32710
$ws.Range("N69").Value = -34208
$ws.Range("H72").Value = 32710
$ws.Range("J72").Value = 32710
$ws.Range("L72").Value = 98130
$ws.Range("N72").Value = -105618
$ws.Range("H107").Value = 608.2
$ws.Range("I107").Value = 585.25
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1755.75
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 164.25
$ws.Range("N107").Value = -5940
